$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new DAMSLTag (col I) and DialogAct (col J)
$updates = @(
    @{Row=8;   I="sd"; J="Statement-non-opinion"},
    @{Row=9;   I="sv"; J="Statement-opinion"},
    @{Row=12;  I="ba"; J="Appreciation"},
    @{Row=15;  I="aa"; J="Agree/Accept"},
    @{Row=18;  I="%";  J="Uninterpretable"},
    @{Row=20;  I="b";  J="Acknowledge (Backchannel)"},
    @{Row=34;  I="sd"; J="Statement-non-opinion"},
    @{Row=44;  I="sd"; J="Statement-non-opinion"},
    @{Row=62;  I="sd"; J="Statement-non-opinion"},
    @{Row=70;  I="ba"; J="Appreciation"},
    @{Row=71;  I="ba"; J="Appreciation"},
    @{Row=76;  I="sd"; J="Statement-non-opinion"},
    @{Row=90;  I="sd"; J="Statement-non-opinion"},
    @{Row=110; I="sd"; J="Statement-non-opinion"},
    @{Row=111; I="sd"; J="Statement-non-opinion"},
    @{Row=112; I="sd"; J="Statement-non-opinion"},
    @{Row=131; I="aa"; J="Agree/Accept"},
    @{Row=151; I="sv"; J="Statement-opinion"}
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.I
    $ws.Range("J" + $u.Row).Value = $u.J
}
